$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the second sheet.
$ws2.Name = "Include #0"

# 2. Update the Version value on the Metadata sheet (row 3).
$ws1.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# 3. Update the Date value on the Metadata sheet (row 8).
$ws1.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# 4. Insert a new "Jurisdiction" row after "Contact" (row 10), pushing
#    Description/Purpose/Copyright/Immutable down by one row.
#    Shift rows 14..11 down to 15..12 (bottom-up to avoid clobbering).
$ws1.Cells.Item(15, 1).Value = $ws1.Cells.Item(14, 1).Value2
$ws1.Cells.Item(15, 2).Value = $ws1.Cells.Item(14, 2).Value2

$ws1.Cells.Item(14, 1).Value = $ws1.Cells.Item(13, 1).Value2
$ws1.Cells.Item(14, 2).Value = $ws1.Cells.Item(13, 2).Value2

$ws1.Cells.Item(13, 1).Value = $ws1.Cells.Item(12, 1).Value2
$ws1.Cells.Item(13, 2).Value = $ws1.Cells.Item(12, 2).Value2

$ws1.Cells.Item(12, 1).Value = $ws1.Cells.Item(11, 1).Value2
$ws1.Cells.Item(12, 2).Value = $ws1.Cells.Item(11, 2).Value2

$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""

# Re-apply the standard data-row formatting (copied from row 10, which
# already carries it) to the new row 11 and the newly extended row 15 so
# no new style entries get created and the look stays consistent.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Range("A15:B15").PasteSpecial(-4122)
